$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.245.35'
$ws.Range('E2').Value = '  +0.65%  '

$ws.Range('D3').Value = '3.799.26'
$ws.Range('E3').Value = '  -0.20%  '

$ws.Range('E4').Value = '  +0.30%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '601.12'
$ws.Range('E5').Value = '  +0.63%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '165.33'
$ws.Range('E6').Value = '  -1.27%  '

$ws.Range('E7').Value = '  -0.26%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.518'
$ws.Range('E8').Value = '  -0.77%  '

$ws.Range('E9').Value = '  -1.13%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.451'
$ws.Range('E10').Value = '  +0.27%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.48'
$ws.Range('E11').Value = '  +2.74%  '

$ws.Range('E12').Value = '  -1.69%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '35.83'
$ws.Range('E13').Value = '  -0.96%  '

$ws.Range('D14').Value = '4.434.16'
$ws.Range('E14').Value = '  -0.23%  '

$ws.Range('D15').Value = '3.787.67'
$ws.Range('E15').Value = '  -0.64%  '

$ws.Range('D16').Value = '68.182.96'
$ws.Range('E16').Value = '  +0.60%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.44'
$ws.Range('E17').Value = '  -0.70%  '

$ws.Range('E18').Value = '  +2.28%  '

$ws.Range('E19').Value = '  -0.50%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '461.75'

$ws.Range('E21').Value = '  -2.02%  '

$ws.Range('E22').Value = '  -0.25%  '

$ws.Range('E23').Value = '  -3.79%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.02'
$ws.Range('E24').Value = '  -0.72%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.07'
$ws.Range('E25').Value = '  -0.45%  '

$ws.Range('E26').Value = '  -0.27%  '

$ws.Range('E27').Value = '  -0.19%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.99'
$ws.Range('E28').Value = '  -0.40%  '

$ws.Range('D29').Value = '3.947.66'
$ws.Range('E29').Value = '  -0.09%  '

$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.66'
$ws.Range('E30').Value = '  -4.98%  '

$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.37'
$ws.Range('E31').Value = '  +1.18%  '

$ws.Range('E32').Value = '  -1.59%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '29.38'
$ws.Range('E33').Value = '  -1.34%  '

$ws.Range('E34').Value = '  +0.02%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.01'
$ws.Range('E35').Value = '  -0.97%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0997'
$ws.Range('E36').Value = '  -0.61%  '

$ws.Range('E37').Value = '  +0.43%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.31'
$ws.Range('E38').Value = '  -3.32%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.81'
$ws.Range('E39').Value = '  +0.14%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.988'
$ws.Range('E40').Value = '  -0.87%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.09%  '

$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.300'
$ws.Range('E43').Value = '  +0.47%  '

$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '47.52'
$ws.Range('E44').Value = '  -1.32%  '

$ws.Range('B45').Value = 'Arweave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '43.50'
$ws.Range('E45').Value = '  -1.29%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '151.74'
$ws.Range('E46').Value = '  +0.89%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.37'
$ws.Range('E47').Value = '  +0.29%  '

$ws.Range('E48').Value = '  +2.35%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '396.18'
$ws.Range('E49').Value = '  -0.43%  '

$ws.Range('E50').Value = '  +5.65%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '26.57'
$ws.Range('E51').Value = '  +0.78%  '
